$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10087.353
$ws.Range("I32").Value = 11117.125
$ws.Range("J32").Value = 9172
$ws.Range("K32").Value = 11117.125
$ws.Range("L32").Value = 9172
$ws.Range("M32").Value = -10791.125
$ws.Range("N32").Value = -9824
$ws.Range("H40").Value = 4204.0835
$ws.Range("J40").Value = 5408.1665
$ws.Range("L40").Value = 5408.1665
$ws.Range("N40").Value = -5758.1665
$ws.Range("H86").Value = 2430.5715
$ws.Range("I86").Value = 2556.75
$ws.Range("J86").Value = 2262.3333
$ws.Range("K86").Value = 2556.75
$ws.Range("L86").Value = 2262.3333
$ws.Range("M86").Value = -1433.75
$ws.Range("N86").Value = -4508.3333
$ws.Range("H89").Value = 2430.5715
$ws.Range("I89").Value = 2556.75
$ws.Range("J89").Value = 2262.3333
$ws.Range("K89").Value = 12783.75
$ws.Range("L89").Value = 11311.6665
$ws.Range("M89").Value = -7167.75
$ws.Range("N89").Value = -22543.6665
$ws.Range("H125").Value = 665
$ws.Range("J125").Value = 665
$ws.Range("L125").Value = 5985
$ws.Range("N125").Value = -10905
$ws.Range("H132").Value = 62502532
$ws.Range("I132").Value = 62502532
$ws.Range("K132").Value = 187507596
$ws.Range("M132").Value = -187505066

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2967.125
$ws.Range("I45").Value = 2333.8572
$ws.Range("K45").Value = 2333.8572
$ws.Range("M45").Value = -1956.8572
$ws.Range("H132").Value = 6872.8237
$ws.Range("I132").Value = 6333.6553
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 19000.9659
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -16470.9659
$ws.Range("N132").Value = -35060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3319.8572
$ws.Range("I20").Value = 3377.8
$ws.Range("K20").Value = 3377.8
$ws.Range("M20").Value = -3130.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10979.1
$ws.Range("I86").Value = 4974
$ws.Range("K86").Value = 4974
$ws.Range("M86").Value = -3851
$ws.Range("H89").Value = 10979.1
$ws.Range("I89").Value = 4974
$ws.Range("K89").Value = 24870
$ws.Range("M89").Value = -19254
$ws.Range("H105").Value = 1469.8334
$ws.Range("I105").Value = 923.5
$ws.Range("J105").Value = 2562.5
$ws.Range("K105").Value = 923.5
$ws.Range("L105").Value = 2562.5
$ws.Range("M105").Value = 823.5
$ws.Range("N105").Value = -6056.5
$ws.Range("H106").Value = 25833
$ws.Range("J106").Value = 25833
$ws.Range("L106").Value = 25833
$ws.Range("N106").Value = -28357
$ws.Range("H107").Value = 941.1177
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 1599.6666
$ws.Range("K107").Value = 800
$ws.Range("L107").Value = 1599.6666
$ws.Range("M107").Value = 1120
$ws.Range("N107").Value = -5439.6666
$ws.Range("H116").Value = 55830.918
$ws.Range("J116").Value = 55830.918
$ws.Range("L116").Value = 55830.918
$ws.Range("N116").Value = -65008.918
$ws.Range("H122").Value = 1147.1333
$ws.Range("I122").Value = 761.1
$ws.Range("K122").Value = 2283.3
$ws.Range("M122").Value = 166.6999999999998
$ws.Range("H132").Value = 5585.5557
$ws.Range("I132").Value = 5080.077
$ws.Range("J132").Value = 6899.8
$ws.Range("K132").Value = 15240.231
$ws.Range("L132").Value = 20699.4
$ws.Range("M132").Value = -12710.231
$ws.Range("N132").Value = -25759.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 225
$ws.Range("J22").Value = 225
$ws.Range("L22").Value = 675
$ws.Range("N22").Value = -1013
$ws.Range("H27").Value = 225
$ws.Range("J27").Value = 225
$ws.Range("L27").Value = 675
$ws.Range("N27").Value = -879
$ws.Range("H43").Value = 5000
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15228
$ws.Range("H113").Value = 1119.3518
$ws.Range("I113").Value = 1136.4166
$ws.Range("K113").Value = 3409.2498
$ws.Range("M113").Value = -1239.2498

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 30000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H80").Value = 6149.7
$ws.Range("I80").Value = 5962.25
$ws.Range("J80").Value = 6899.5
$ws.Range("K80").Value = 5962.25
$ws.Range("L80").Value = 6899.5
$ws.Range("M80").Value = -4964.25
$ws.Range("N80").Value = -8895.5
$ws.Range("H83").Value = 6149.7
$ws.Range("I83").Value = 5962.25
$ws.Range("J83").Value = 6899.5
$ws.Range("K83").Value = 29811.25
$ws.Range("L83").Value = 34497.5
$ws.Range("M83").Value = -24819.25
$ws.Range("N83").Value = -44481.5
$ws.Range("H102").Value = 2392.182
$ws.Range("I102").Value = 1835.2222
$ws.Range("K102").Value = 1835.2222
$ws.Range("M102").Value = -213.2221999999999
$ws.Range("H132").Value = 3750.35
$ws.Range("I132").Value = 3706.4119
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 11119.2357
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -8589.235700000001
$ws.Range("N132").Value = -17057.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 827.7
$ws.Range("I16").Value = 847.44446
$ws.Range("J16").Value = 650
$ws.Range("K16").Value = 847.44446
$ws.Range("L16").Value = 650
$ws.Range("M16").Value = -677.44446
$ws.Range("N16").Value = -990
$ws.Range("H122").Value = 3237.5
$ws.Range("I122").Value = 3237.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9712.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7262.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2551.963
$ws.Range("I132").Value = 2294.3635
$ws.Range("K132").Value = 6883.0905
$ws.Range("M132").Value = -4353.0905

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50500000
$ws.Range("H113").Value = 547.3333
$ws.Range("I113").Value = 422.33334
$ws.Range("J113").Value = 797.3333
$ws.Range("K113").Value = 1267.00002
$ws.Range("L113").Value = 2391.9999
$ws.Range("M113").Value = 902.9999800000001
$ws.Range("N113").Value = -6731.9999

Write-Output "done"